$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The author moved the active cell from E7 to E6 before saving.
$ws.Range("E6").Select()

# Row 6 data was corrected: fewer remaining cases (10 -> 7) and the
# status moved from "Suited to Manual" to "Outdated".
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = "Outdated"
